{"js": "// \"Remove another reference to all-caps titles and authors in the Word\n// template\" -- the Guidelines paragraph no longer claims the title/authors\n// must be in all capitals:\n//   \"The title should be all capitals, bold, ...\"\n//     -> \"The title should be capitalized like a sentence, bold, ...\"\n//   \"The list of authors should be all capitals 9 point font, ...\"\n//     -> \"The list of authors should be 9 point font, ...\"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldText =\n  \"We have provided a file showing the layout required when producing a \" +\n  \"paper using Word (this is likely to be the document you are currently \" +\n  \"reading). This is the body text and should be 10pt font. The title \" +\n  \"should be all capitals, bold, 11pt font. The list of authors should be \" +\n  \"all capitals 9 point font, and the list of affiliations should be \" +\n  \"regular 9 point font. You can delete our sample text and replace it \" +\n  \"with your own contribution to the volume, although we recommend that \" +\n  \"you keep an initial version of this file for reference. \";\n\nconst newText =\n  \"We have provided a file showing the layout required when producing a \" +\n  \"paper using Word (this is likely to be the document you are currently \" +\n  \"reading). This is the body text and should be 10pt font. The title \" +\n  \"should be capitalized like a sentence, bold, 11pt font. The list of \" +\n  \"authors should be 9 point font, and the list of affiliations should be \" +\n  \"regular 9 point font. You can delete our sample text and replace it \" +\n  \"with your own contribution to the volume, although we recommend that \" +\n  \"you keep an initial version of this file for reference. \";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === oldText) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  // Fall back to a loose match in case whitespace differs slightly.\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text.indexOf(\"We have provided a file\") !== -1) {\n      target = paragraphs.items[i];\n      break;\n    }\n  }\n}\n\nif (target) {\n  target.insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# \"Remove another reference to all-caps titles and authors in the Word\n# template\" -- the Guidelines paragraph no longer claims the title/authors\n# must be in all capitals:\n#   \"The title should be all capitals, bold, ...\"\n#     -> \"The title should be capitalized like a sentence, bold, ...\"\n#   \"The list of authors should be all capitals 9 point font, ...\"\n#     -> \"The list of authors should be 9 point font, ...\"\n$d = $word.ActiveDocument\n\n$find1 = $d.Content.Find\n$find1.Text = \"all capitals, bold\"\n$find1.Replacement.Text = \"capitalized like a sentence, bold\"\n$find1.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n$find2 = $d.Content.Find\n$find2.Text = \"The list of authors should be all capitals 9 point font\"\n$find2.Replacement.Text = \"The list of authors should be 9 point font\"\n$find2.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n"}
